# Apply the edits described in the commit:
# "more edits - need to figure out problems with which cobbles should be
#  included in the analysis for E2 and possibly E3"

$wb = $excel.ActiveWorkbook

# 1) Workbook-level compatibility flag (workbookPr checkCompatibility="1")
$wb.CheckCompatibility = $true

# 2) Rename sheet "R input" -> "During experiment recording_fix"
$ws = $wb.Worksheets.Item("R input")
$ws.Name = "During experiment recording_fix"

# Make sure this sheet is the active one (it was tabSelected before/after).
$ws.Activate()

# 3) Cell edits on the renamed sheet.

# G113 used to hold the text "200+" from the shared-string table; it now
# holds the plain number 200.
$ws.Range("G113").Value = 200

# Row 134
$ws.Range("E134").Value = 115.78
$ws.Range("G134").Value = 59

# Row 135
$ws.Range("G135").Value = 281

# Row 136
$ws.Range("G136").Value = 299

# Row 137 (also shrinks the shared formula range for I137 down to just I137)
$ws.Range("E137").Value = 76.36
$ws.Range("G137").Value = 106
$ws.Range("I137").Formula = "=1/10"

# Row 138
$ws.Range("E138").Value = 116.47
$ws.Range("G138").Value = 138

# Row 139
$ws.Range("G139").Value = 4

# Row 140 (also shrinks the shared formula range for I140 down to just I140)
$ws.Range("G140").Value = 262
$ws.Range("I140").Formula = "=1/10"

# Row 141
$ws.Range("G141").Value = 320

# Row 142
$ws.Range("G142").Value = 197

# Row 143
$ws.Range("G143").Value = 292

# 4) Update the sheet's view/selection: move the selection to G136 (this
#    also drops the stale topLeftCell="A165" scroll position).
$ws.Range("G136").Select()

Write-Output "edits applied"
